# Updated symbol list (coin price/volume snapshot refresh) on
# Fri Feb 17 11:53:06 UTC 2023 with GitHub Actions.
#
# Price (D) and Volume(1h) (E) are stored as literal text in this sheet
# (e.g. "309.87", "-3.17%"), not as numbers/percentages. Writing a plain
# numeric-looking string via .Value would make Excel auto-coerce it to a
# real number, so we prefix with a quote (like typing '309.87 into the
# UI) to force text, then ClearFormats() to drop the resulting quote-prefix
# style so the cell keeps its original (unstyled) formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'309.87"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'-3.17%"
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'50.77"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'3.85%"
$ws.Range("E3").ClearFormats()
$ws.Range("D4").Value = "'5.160"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "'-1.73%"
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'0.07780"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'-3.76%"
$ws.Range("E5").ClearFormats()
$ws.Range("E6").Value = "'-2.10%"
$ws.Range("E6").ClearFormats()
$ws.Range("D7").Value = "'1.350"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'12.00%"
$ws.Range("E7").ClearFormats()
$ws.Range("D8").Value = "'1.567"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'-4.71%"
$ws.Range("E8").ClearFormats()
$ws.Range("D9").Value = "'0.1210"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'-6.24%"
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = "'0.1984"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'2.68%"
$ws.Range("E10").ClearFormats()
$ws.Range("D11").Value = "'0.09604"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'2.57%"
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = "'0.04740"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'3.26%"
$ws.Range("E12").ClearFormats()
$ws.Range("E13").Value = "'-0.50%"
$ws.Range("E13").ClearFormats()
$ws.Range("E14").Value = "'-4.80%"
$ws.Range("E14").ClearFormats()
$ws.Range("D15").Value = "'0.005787"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'-2.49%"
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = "'0.007487"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'2,012.26%"
$ws.Range("E16").ClearFormats()
$ws.Range("E17").Value = "'-0.36%"
$ws.Range("E17").ClearFormats()
$ws.Range("E18").Value = "'0.47%"
$ws.Range("E18").ClearFormats()
$ws.Range("E19").Value = "'1.98%"
$ws.Range("E19").ClearFormats()
$ws.Range("D20").Value = "'8.008"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'-1.15%"
$ws.Range("E20").ClearFormats()
$ws.Range("D21").Value = "'0.1372"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'-0.90%"
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = "'0.3094"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'-1.01%"
$ws.Range("E22").ClearFormats()
$ws.Range("D23").Value = "'0.04164"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'-0.11%"
$ws.Range("E23").ClearFormats()
$ws.Range("E24").Value = "'-2.67%"
$ws.Range("E24").ClearFormats()
$ws.Range("E25").Value = "'-6.88%"
$ws.Range("E25").ClearFormats()
$ws.Range("D26").Value = "'0.0001350"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "'-0.19%"
$ws.Range("E26").ClearFormats()
$ws.Range("D38").Value = "'0.02595"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'-3.87%"
$ws.Range("E38").ClearFormats()
$ws.Range("D39").Value = "'0.06013"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'6.09%"
$ws.Range("E39").ClearFormats()
$ws.Range("D40").Value = "'0.01100"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'74.11%"
$ws.Range("E40").ClearFormats()
$ws.Range("D41").Value = "'0.007858"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'-0.73%"
$ws.Range("E41").ClearFormats()
$ws.Range("E42").Value = "'-1.25%"
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = "'0.008384"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'8.64%"
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = "'0.007676"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'-5.25%"
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'0.3389"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'6.10%"
$ws.Range("E45").ClearFormats()
$ws.Range("D46").Value = "'0.00007360"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'6.51%"
$ws.Range("E46").ClearFormats()
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'-0.19%"
$ws.Range("E47").ClearFormats()
$ws.Range("B48").Value = "CoinbaseStockToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D48").Value = "'0.002619"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'-34.63%"
$ws.Range("E48").ClearFormats()
$ws.Range("B49").Value = "BOLO"
$ws.Range("C49").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D49").Value = "'0.05318"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'-7.21%"
$ws.Range("E49").ClearFormats()
$ws.Range("D50").Value = "'0.00002099"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'-0.19%"
$ws.Range("E50").ClearFormats()
$ws.Range("E51").Value = "'-0.19%"
$ws.Range("E51").ClearFormats()
